# Scoreboard.xlsx edit: add Semi-Final sheets (FM / FF) + "First Stage Points" /
# "Semi Final Points" columns to the existing SFM / SFF sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) SFM ("SFM" tab) - insert a "First Stage Points" column (new col B,
#    Snatch/Clean&Jerk shift right to C/D) and refresh the qualifier table
#    for the men's semi-final (different / reordered set of teams).
# ---------------------------------------------------------------------------
$sfm = $wb.Worksheets.Item("SFM")
$sfm.Columns("B").Insert()
$sfm.Range("B1").Value = "First Stage Points"

$sfmData = @(
  @("Håkon Konningen og Njål Christensen", 8, 75, 100),
  @("Anders Magnus Nes og Anders Vinnes Jacobsen", 9, 78, 111),
  @("Anders J. Svalestuen og Gabriel Kristiansen", 9, 80, 99),
  @("Ole Andre Elvebakk og Georg Kongsvik", 11, 66, 85),
  @("Magnus Ødegaard og Kornelius Skrettingland", 8, 52, 98),
  @("Vegard Austrheim Vågen og Henrik Eliassen", 8, 77, 104)
)

$r = 2
foreach ($row in $sfmData) {
  $sfm.Range("A$r").Value = $row[0]
  $sfm.Range("B$r").Value = $row[1]
  $sfm.Range("C$r").Value = $row[2]
  $sfm.Range("D$r").Value = $row[3]
  $r++
}

$sfm.Columns("B").ColumnWidth = 16.14
$sfm.PageSetup.PaperSize = 9
$sfm.PageSetup.Orientation = 1
$sfm.Range("D18").Select()

# ---------------------------------------------------------------------------
# 2) SFF ("SFF" tab) - same column insert, and the qualifier table shrinks
#    from 7 teams to 6.
# ---------------------------------------------------------------------------
$sff = $wb.Worksheets.Item("SFF")
$sff.Columns("B").Insert()
$sff.Range("B1").Value = "First Stage Points"
$sff.Rows("8").Delete()

$sffData = @(
  @("Beata Wilman og Ingrid Hamnes", 8, 66, 90),
  @("Renate Berntsen Hansen og Karoline Granås", 8, 59, 85),
  @("Maria Hanssen og Cecilie Rabben", 9, 72, 81),
  @("Marianne U. Henriksen og Mari S. Andersen", 9, 78, 75),
  @("Sara Yuzer og Martine Baalsrud", 9, 55, 91),
  @("Frid Kaspersen og Renate Loraas", 11, 69, 80)
)

$r = 2
foreach ($row in $sffData) {
  $sff.Range("A$r").Value = $row[0]
  $sff.Range("B$r").Value = $row[1]
  $sff.Range("C$r").Value = $row[2]
  $sff.Range("D$r").Value = $row[3]
  $r++
}

$sff.Columns("B").ColumnWidth = 16.14
$sff.Range("B21").Select()

# ---------------------------------------------------------------------------
# 3) Two new sheets for the semi-final round, inserted right after "SFF" and
#    before "ScoreMatrix": "FM" (men) and "FF" (women).
# ---------------------------------------------------------------------------
$fm = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sff)
$fm.Name = "FM"
$fm.Range("A1").Value = "Team"
$fm.Range("B1").Value = "First Stage Points"
$fm.Range("C1").Value = "Semi Final Points"
$fm.Range("D1").Value = "Minute1"
$fm.Range("E1").Value = "Second1"
$fm.Range("F1").Value = "Rep1"
$fm.Columns("A").ColumnWidth = 5.85546875
$fm.Columns("B").ColumnWidth = 16.28515625
$fm.Columns("C").ColumnWidth = 16.28515625
$fm.Columns("D").ColumnWidth = 8.42578125
$fm.Columns("E").ColumnWidth = 8.42578125
$fm.Columns("F").ColumnWidth = 5.42578125
$fm.Range("E20").Select()

$ff = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $fm)
$ff.Name = "FF"
$ff.Range("A1").Value = "Team"
$ff.Range("B1").Value = "First Stage Points"
$ff.Range("C1").Value = "Semi Final Points"
$ff.Range("D1").Value = "Minute1"
$ff.Range("E1").Value = "Second1"
$ff.Range("F1").Value = "Rep1"
$ff.Columns("A").ColumnWidth = 5.85546875
$ff.Columns("B").ColumnWidth = 16.28515625
$ff.Columns("C").ColumnWidth = 16.28515625
$ff.Columns("D").ColumnWidth = 8.42578125
$ff.Columns("E").ColumnWidth = 8.42578125
$ff.Columns("F").ColumnWidth = 5.42578125
$ff.Cells.Select()

# ---------------------------------------------------------------------------
# 4) "ScoreMatrix" tab keeps its data; only the remembered selection moves.
# ---------------------------------------------------------------------------
$scoreMatrix = $wb.Worksheets.Item("ScoreMatrix")
$scoreMatrix.Range("F19").Select()

# ---------------------------------------------------------------------------
# 5) View bookkeeping on the two original tabs (ScoreM / ScoreF) and which
#    tab is active when the workbook is reopened (SFF, per activeTab=3).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("ScoreM").Range("C1:E1").Select()
$wb.Worksheets.Item("ScoreF").Range("C19").Select()

$sff.Activate()

Write-Output "edit complete"
